$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value  = 8.350541
$ws.Cells.Item(2, 8).Value  = 25.051623
$ws.Cells.Item(2, 9).Value  = 0.3086892463293835
$ws.Cells.Item(2, 10).Value = 0.3086892463293835
$ws.Cells.Item(2, 15).Value = 0.7091726973716084
$ws.Cells.Item(2, 16).Value = 0.7091726973716084
$ws.Cells.Item(2, 17).Value = 14.26094257925333
$ws.Cells.Item(2, 18).Value = 128.34848321328
$ws.Cells.Item(2, 19).Value = 0.2189139854690178
$ws.Cells.Item(2, 20).Value = 0.2189139854690178

# Row 3
$ws.Cells.Item(3, 7).Value  = 8.350541
$ws.Cells.Item(3, 8).Value  = 25.051623
$ws.Cells.Item(3, 9).Value  = 0.3086892463293835
$ws.Cells.Item(3, 10).Value = 0.3086892463293835
$ws.Cells.Item(3, 13).Value = 0.7003526666666667
$ws.Cells.Item(3, 14).Value = 2.101058
$ws.Cells.Item(3, 15).Value = 0.2908273026283917
$ws.Cells.Item(3, 16).Value = 0.2908273026283917
$ws.Cells.Item(3, 17).Value = 5.848323657459334
$ws.Cells.Item(3, 18).Value = 52.634912917134
$ws.Cells.Item(3, 19).Value = 0.08977526086036577
$ws.Cells.Item(3, 20).Value = 0.08977526086036577

# Row 4
$ws.Cells.Item(4, 9).Value  = 0.4146406124520329
$ws.Cells.Item(4, 10).Value = 0.4146406124520329
$ws.Cells.Item(4, 15).Value = 0.7091726973716084
$ws.Cells.Item(4, 16).Value = 0.7091726973716084
$ws.Cells.Item(4, 19).Value = 0.2940518015724239
$ws.Cells.Item(4, 20).Value = 0.2940518015724239

# Row 5
$ws.Cells.Item(5, 9).Value  = 0.4146406124520329
$ws.Cells.Item(5, 10).Value = 0.4146406124520329
$ws.Cells.Item(5, 13).Value = 0.7003526666666667
$ws.Cells.Item(5, 14).Value = 2.101058
$ws.Cells.Item(5, 15).Value = 0.2908273026283917
$ws.Cells.Item(5, 16).Value = 0.2908273026283917
$ws.Cells.Item(5, 17).Value = 7.855642954789333
$ws.Cells.Item(5, 18).Value = 70.700786593104
$ws.Cells.Item(5, 19).Value = 0.120588810879609
$ws.Cells.Item(5, 20).Value = 0.120588810879609

# Row 6
$ws.Cells.Item(6, 7).Value  = 7.484372666666666
$ws.Cells.Item(6, 8).Value  = 22.453118
$ws.Cells.Item(6, 9).Value  = 0.2766701412185836
$ws.Cells.Item(6, 10).Value = 0.2766701412185836
$ws.Cells.Item(6, 15).Value = 0.7091726973716084
$ws.Cells.Item(6, 16).Value = 0.7091726973716084
$ws.Cells.Item(6, 17).Value = 12.78171184849778
$ws.Cells.Item(6, 18).Value = 115.03540663648
$ws.Cells.Item(6, 19).Value = 0.1962069103301667
$ws.Cells.Item(6, 20).Value = 0.1962069103301667

# Row 7
$ws.Cells.Item(7, 7).Value  = 7.484372666666666
$ws.Cells.Item(7, 8).Value  = 22.453118
$ws.Cells.Item(7, 9).Value  = 0.2766701412185836
$ws.Cells.Item(7, 10).Value = 0.2766701412185836
$ws.Cells.Item(7, 13).Value = 0.7003526666666667
$ws.Cells.Item(7, 14).Value = 2.101058
$ws.Cells.Item(7, 15).Value = 0.2908273026283917
$ws.Cells.Item(7, 16).Value = 0.2908273026283917
$ws.Cells.Item(7, 17).Value = 5.241700355427112
$ws.Cells.Item(7, 18).Value = 47.175303198844
$ws.Cells.Item(7, 19).Value = 0.08046323088841686
$ws.Cells.Item(7, 20).Value = 0.08046323088841686
